$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 399.0625
$ws.Range("I28").Value = 392.33334
$ws.Range("J28").Value = 500
$ws.Range("K28").Value = 392.33334
$ws.Range("L28").Value = 500
$ws.Range("M28").Value = 92.66665999999998
$ws.Range("N28").Value = -1470
$ws.Range("H62").Value = 2667.4614
$ws.Range("I62").Value = 1966.625
$ws.Range("J62").Value = 3788.8
$ws.Range("K62").Value = 1966.625
$ws.Range("L62").Value = 3788.8
$ws.Range("M62").Value = -1342.625
$ws.Range("N62").Value = -5036.8
$ws.Range("H65").Value = 2667.4614
$ws.Range("I65").Value = 1966.625
$ws.Range("J65").Value = 3788.8
$ws.Range("K65").Value = 9833.125
$ws.Range("L65").Value = 18944
$ws.Range("M65").Value = -6713.125
$ws.Range("N65").Value = -25184
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()
$ws.Range("H88").Value = 8386.6875
$ws.Range("I88").Value = 6963.2856
$ws.Range("J88").Value = 9493.777
$ws.Range("K88").Value = 6963.2856
$ws.Range("L88").Value = 9493.777
$ws.Range("M88").Value = -6557.2856
$ws.Range("N88").Value = -10305.777
$ws.Range("H91").Value = 8386.6875
$ws.Range("I91").Value = 6963.2856
$ws.Range("J91").Value = 9493.777
$ws.Range("K91").Value = 6963.2856
$ws.Range("L91").Value = 9493.777
$ws.Range("M91").Value = -5559.2856
$ws.Range("N91").Value = -12301.777
$ws.Range("H111").Value = 1920.2222
$ws.Range("I111").Value = 1847.75
$ws.Range("J111").Value = 2500
$ws.Range("K111").Value = 5543.25
$ws.Range("L111").Value = 7500
$ws.Range("M111").Value = -2476.25
$ws.Range("N111").Value = -13634
$ws.Range("H113").Value = 5067.4375
$ws.Range("I113").Value = 4959.875
$ws.Range("J113").Value = 5175
$ws.Range("K113").Value = 4959.875
$ws.Range("L113").Value = 5175
$ws.Range("M113").Value = -1705.875
$ws.Range("N113").Value = -11683
$ws.Range("H141").Value = 535948.1
$ws.Range("I141").Value = 1461.6316
$ws.Range("J141").Value = 1986697.2
$ws.Range("K141").Value = 4384.8948
$ws.Range("L141").Value = 5960091.6
$ws.Range("M141").Value = 795.1052
$ws.Range("N141").Value = -5970451.6
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4512.02
$ws.Range("I32").Value = 3394.9167
$ws.Range("J32").Value = 10767.8
$ws.Range("K32").Value = 3394.9167
$ws.Range("L32").Value = 10767.8
$ws.Range("M32").Value = -3107.9167
$ws.Range("N32").Value = -11341.8
$ws.Range("H45").Value = 1632
$ws.Range("I45").Value = 1068.4138
$ws.Range("J45").Value = 3675
$ws.Range("K45").Value = 1068.4138
$ws.Range("L45").Value = 3675
$ws.Range("M45").Value = -691.4138
$ws.Range("N45").Value = -4429
$ws.Range("H110").Value = 1562.409
$ws.Range("I110").Value = 660
$ws.Range("J110").Value = 3968.8333
$ws.Range("K110").Value = 660
$ws.Range("L110").Value = 3968.8333
$ws.Range("M110").Value = 1385
$ws.Range("N110").Value = -8058.8333
$ws.Range("H118").Value = 29400
$ws.Range("J118").Value = 29400
$ws.Range("L118").Value = 29400
$ws.Range("N118").Value = -32714
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 6026.2666
$ws.Range("I20").Value = 8647
$ws.Range("J20").Value = 3031.1428
$ws.Range("K20").Value = 8647
$ws.Range("L20").Value = 3031.1428
$ws.Range("M20").Value = -8400
$ws.Range("N20").Value = -3525.1428
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 4522.5
$ws.Range("I62").Value = 2860
$ws.Range("J62").Value = 5520
$ws.Range("K62").Value = 2860
$ws.Range("L62").Value = 5520
$ws.Range("M62").Value = -2236
$ws.Range("N62").Value = -6768
$ws.Range("H65").Value = 4522.5
$ws.Range("I65").Value = 2860
$ws.Range("J65").Value = 5520
$ws.Range("K65").Value = 14300
$ws.Range("L65").Value = 27600
$ws.Range("M65").Value = -11180
$ws.Range("N65").Value = -33840
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 341.83334
$ws.Range("J17").Value = 673
$ws.Range("L17").Value = 2019
$ws.Range("N17").Value = -2357
$ws.Range("H34").Value = 8559.23
$ws.Range("I34").Value = 210
$ws.Range("J34").Value = 18300
$ws.Range("K34").Value = 630
$ws.Range("L34").Value = 54900
$ws.Range("M34").Value = -546
$ws.Range("N34").Value = -55068
$ws.Range("H39").Value = 1101.3334
$ws.Range("I39").Value = 300
$ws.Range("J39").Value = 1502
$ws.Range("K39").Value = 900
$ws.Range("L39").Value = 4506
$ws.Range("M39").Value = -606
$ws.Range("N39").Value = -5094
$ws.Range("H55").Value = 1630
$ws.Range("J55").Value = 3940
$ws.Range("L55").Value = 11820
$ws.Range("N55").Value = -12174
$ws.Range("H107").Value = 1306.9
$ws.Range("I107").Value = 782
$ws.Range("J107").Value = 1831.8
$ws.Range("K107").Value = 2346
$ws.Range("L107").Value = 5495.4
$ws.Range("M107").Value = -426
$ws.Range("N107").Value = -9335.4
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 5000
$ws.Range("I52").Value = 0
$ws.Range("K52").Value = 0
$ws.Range("M52").ClearContents()
$ws.Range("H113").Value = 5750
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 5750
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 5750
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -10090
$ws.Range("H122").Value = 3374.2144
$ws.Range("I122").Value = 2814.6316
$ws.Range("J122").Value = 4555.5557
$ws.Range("K122").Value = 8443.8948
$ws.Range("L122").Value = 13666.6671
$ws.Range("M122").Value = -5993.8948
$ws.Range("N122").Value = -18566.6671
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2187.5
$ws.Range("I40").Value = 935
$ws.Range("J40").Value = 3440
$ws.Range("K40").Value = 935
$ws.Range("L40").Value = 3440
$ws.Range("M40").Value = -799
$ws.Range("N40").Value = -3712
$ws.Range("H45").Value = 9000
$ws.Range("I45").Value = 0
$ws.Range("K45").Value = 0
$ws.Range("M45").ClearContents()
$ws.Range("H46").Value = 2309.0908
$ws.Range("I46").Value = 550
$ws.Range("K46").Value = 550
$ws.Range("M46").Value = -362
$ws.Range("H48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("N48").ClearContents()
$ws.Range("H55").Value = 1070.8572
$ws.Range("I55").Value = 270
$ws.Range("J55").Value = 1391.2
$ws.Range("K55").Value = 270
$ws.Range("L55").Value = 1391.2
$ws.Range("M55").Value = -97
$ws.Range("N55").Value = -1737.2
$ws.Range("H61").Value = 100003840
$ws.Range("I61").Value = 166669230
$ws.Range("K61").Value = 166669230
$ws.Range("M61").Value = -166669028
$ws.Range("H113").Value = 100003840
$ws.Range("I113").Value = 166669230
$ws.Range("K113").Value = 166669230
$ws.Range("M113").Value = -166667060
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H6").Value = 11857.429
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 11857.429
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 11857.429
$ws.Range("M6").ClearContents()
$ws.Range("N6").Value = -12087.429
$ws.Range("H62").Value = 2920
$ws.Range("I62").Value = 2920
$ws.Range("K62").Value = 2920
$ws.Range("M62").Value = -2296
$ws.Range("H65").Value = 2920
$ws.Range("I65").Value = 2920
$ws.Range("K65").Value = 14600
$ws.Range("M65").Value = -11480
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()
